$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking price string to a cell while keeping it
# stored as TEXT (matching the source sheet, where every "Price" cell is an
# inline/shared string, not a number). Temporarily flip the cell to the
# Text number format so Excel doesn't auto-coerce the string into a double,
# then restore the style so no stray formatting sticks around.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Simple price (Price column, "D") updates
Set-TextValue "D2"  "236.20"
Set-TextValue "D3"  "21.78"
Set-TextValue "D4"  "5.354"
Set-TextValue "D6"  "6.462"
Set-TextValue "D7"  "3.352"
Set-TextValue "D8"  "0.7987"
Set-TextValue "D9"  "1.038"
Set-TextValue "D10" "0.1392"
Set-TextValue "D11" "0.07317"
Set-TextValue "D12" "0.03188"
Set-TextValue "D13" "0.02965"
Set-TextValue "D14" "0.09238"
Set-TextValue "D15" "0.001660"
Set-TextValue "D16" "3.257"
Set-TextValue "D17" "0.04788"
Set-TextValue "D18" "0.0005714"
Set-TextValue "D19" "0.006219"
Set-TextValue "D20" "0.005043"
Set-TextValue "D21" "0.001050"
Set-TextValue "D22" "0.0001501"
Set-TextValue "D23" "0.0003996"
Set-TextValue "D24" "3.948"
Set-TextValue "D27" "0.1045"
Set-TextValue "D40" "0.04115"
Set-TextValue "D41" "0.006996"
Set-TextValue "D44" "0.008818"
Set-TextValue "D45" "0.00005436"
Set-TextValue "D46" "0.00000000750"
Set-TextValue "D47" "0.6756"
Set-TextValue "D48" "0.03473"
Set-TextValue "D49" "0.00002101"

# Rows 42 and 43 swap coins: CEJI <-> BKEXToken (name, link, rank-prefixed
# volume label), each with its own refreshed price.
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1038"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.002953"
$ws.Range("E43").Value = "42CEJICEJI"
